$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# =====================================================================
# Sheet1 "Current expenditure": insert 5 new region rows (Zanzibar
# regions), each with a blank spending value, keeping alphabetical
# ordering of the region list.
# =====================================================================

# Insert rows for "Kaskazini Pemba" (row5) and "Kaskazini Unguja" (row6)
# before current row 5 (Katavi). Values are written Unguja-then-Pemba so
# the shared-string table picks up the same allocation order as the
# original edit (idx 39 = Unguja, idx 40 = Pemba) even though Pemba ends
# up sorted above Unguja in the sheet.
$ws1.Rows.Item(5).Insert()
$ws1.Range("A6:B6").Copy()
$ws1.Range("A5:B5").PasteSpecial(-4122)
$ws1.Range("E6").Copy()
$ws1.Range("E5").PasteSpecial(-4122)

$ws1.Rows.Item(6).Insert()
$ws1.Range("A7:B7").Copy()
$ws1.Range("A6:B6").PasteSpecial(-4122)
$ws1.Range("E7").Copy()
$ws1.Range("E6").PasteSpecial(-4122)

$ws1.Range("A6").Value2 = "Kaskazini Unguja"
$ws1.Range("B6").Value2 = $null
$ws1.Range("A5").Value2 = "Kaskazini Pemba"
$ws1.Range("B5").Value2 = $null

# Insert "Kusini Pemba" before current row 9 (Lindi, now shifted)
$ws1.Rows.Item(9).Insert()
$ws1.Range("A10:B10").Copy()
$ws1.Range("A9:B9").PasteSpecial(-4122)
$ws1.Range("E10").Copy()
$ws1.Range("E9").PasteSpecial(-4122)
$ws1.Range("A9").Value2 = "Kusini Pemba"
$ws1.Range("B9").Value2 = $null

# Insert "Kusini Unguja" before current row 10 (Lindi, now shifted again)
$ws1.Rows.Item(10).Insert()
$ws1.Range("A11:B11").Copy()
$ws1.Range("A10:B10").PasteSpecial(-4122)
$ws1.Range("E11").Copy()
$ws1.Range("E10").PasteSpecial(-4122)
$ws1.Range("A10").Value2 = "Kusini Unguja"
$ws1.Range("B10").Value2 = $null

# Insert "Mjini Magharibi" before current row 14 (Morogoro, now shifted)
$ws1.Rows.Item(14).Insert()
$ws1.Range("A15:B15").Copy()
$ws1.Range("A14:B14").PasteSpecial(-4122)
$ws1.Range("E15").Copy()
$ws1.Range("E14").PasteSpecial(-4122)
$ws1.Range("A14").Value2 = "Mjini Magharibi"
$ws1.Range("B14").Value2 = $null

# Extra (non-contiguous) blank styled row below the table
$ws1.Range("B25").Copy()
$ws1.Range("B28").PasteSpecial(-4122)
$ws1.Range("B28").Value2 = $null

# Column A is a touch wider now that longer region names are present
$ws1.Columns.Item(1).ColumnWidth = 14.3

# =====================================================================
# Sheet2 "Optimal funding scenario": clear the "x" marker from the
# "Check to include in analysis" column for the first/second/fourth
# scenario rows (rows 2, 3, 5) - only row 4 keeps it.
# =====================================================================
$ws2.Range("E2").Value2 = $null
$ws2.Range("E3").Value2 = $null
$ws2.Range("E5").Value2 = $null

# =====================================================================
# View state: sheet2 loses tabSelected / its selection moves to D4;
# sheet1 becomes the active tab with a 125% zoom and F5 selected.
# =====================================================================
$ws2.Select()
$ws2.Range("D4").Select()
$ws1.Select()
$ws1.Range("F5").Select()
$excel.ActiveWindow.Zoom = 125

Write-Host "done"
